# Auto-generated edit script applying the curso_elton.xlsx commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Fill in sim/nao + mensagem (column G) for the newly documented lesson rows ---
$ws.Range("C50").Value = "sim"
$ws.Range("D50").Value = "não"
$ws.Range("G50").Value = 'Perfeito nada mais é do que o simples bem feito. Marcelo Singulani'

$ws.Range("C51").Value = "sim"
$ws.Range("D51").Value = "não"
$ws.Range("G51").Value = 'Sua familia vai te apoiar depois de entender o que você faz. Comece explicando e vendendo o que você faz para quem te ama.Helen Caetano'

$ws.Range("C53").Value = "sim"
$ws.Range("D53").Value = "não"
$ws.Range("G53").Value = 'É melhor as pessoas não gostarem de você pelo que você é do que gostarem pelo que você não é. Victor Damasio'

$ws.Range("C55").Value = "sim"
$ws.Range("D55").Value = "não"
$ws.Range("G55").Value = 'É melhor as pessoas não gostarem de você pelo que você é do que gostarem pelo que você não é. Victor Damasio'

$ws.Range("C57").Value = "sim"
$ws.Range("D57").Value = "não"
$ws.Range("G57").Value = 'Por dentro os heróis são como você e eu, pessoas comuns tentando fazer o seu melhor'

$ws.Range("C66").Value = "sim"
$ws.Range("D66").Value = "não"
$ws.Range("G66").Value = 'O unico atalho que existe é a persistência. Paulo Pereira'

# --- Add the new row 70 grand-total line ---
$ws.Range("A70").Value = "total curso"
$ws.Range("B70").Formula = "=SUM(B12,B20,B28,B37,B48,B58,B68)"
$ws.Range("B70").NumberFormat = $ws.Range("B68").NumberFormat()

# --- Add the new reviewer comments ---
$ws.Range("F50").AddComment('Alan Jose Nascimento:
Como construir 
Perfeito nada mais é que o simples perfeito (marcelo singulani)
comunicacao leve
repetição leva perfeição: faça o simples varias vezes 
')
$ws.Range("F51").AddComment('Alan Jose Nascimento:
Perfis viva segunda...
Shiny Happy People: dica --> entenda que a sua mensagem tambem tem que ser importante para o outro e alem de ser verdade ela precisa entendida.
ama aquilo mas nao consegue viver daquilo.
Seraci dica --> tenha coragem de focar e decidir quem é a pessoal que você vai ajudar
toca varias coisas ao mesmo tempo problema de foco, nao tem coragem de identificar um nicho e trabalhar com aquilo. 
Trabalhador sem alma dica dica --> precisa aprender a se importar mais com ele para depois ele dar conta de se importar com o outro
nao consegue se importarm com outro
garimpeiro do proposito: dica --> voce está nessa busca o seu cliente nao, ele quer algo concreto
Super dica: Cuidado para nao colocar a busca do seu proposito na frente do seu cliente
tende a subjetivar as coisas relatizar tudo....
nao tente encontrar seu proposito na sua menasgem...
Cara da Grana: dica --> Tente encontrar um proposito para colocar toda a sua expertise em pratica isso vai te ajudar a se importar com as pessoas, voce vai se conectar com elas , você vai se relacionar com elas e voce nao vai perder nada pelo contrario voce vai ganhar mais dinheiro ainda.
Super dica: quer ganhar mais dinheiro coloque empatia em pratica
Mentalidade focada em resultado, entende quem manda é o cliente, mas ele nao consegue se importar com o cliente.
Nao cria vinculo no relacionamento, e gera comunicação mercenaria
Super Realizado dica -->  Não desista de levar essa mensagem desse  jeito, de entregar o seu produto dessa forma por oque o Super Realizado tem a capacidade de ser Honesto, integro e gerar valor para as pessoas como poucos empreendedores tem 
ele faz o que gosta.!!!
esta ')
$ws.Range("F53").AddComment('Alan Jose Nascimento:
o cliente deve concordar com a gente de pouquinho em pouquinho:
Construa uma ponte entreo seu objetivo e o objetivo do seu cliente.
Insight da aula : 5:04
Não é por que você toca na dor, que você sera a pessoa mais indicada para ajudar a resolver essa minha dor.
Temos que construir uma percepção na cabeça do cliente, por isso devemos montar na pedaço por pedaço.
As pessoas compram um futuro melhor!!1
Temos que conseguir a atençao das pessoas para gerar uma conexão intença e fazer a cabeça balancar a cabeça 3 vezes e se possivel colocar a mão no queixo.....
Exemplo do médico com dor no ombro 
explicar 3 coisas basicas:
! - por que eu cai no buraco;
2 -Por que eu nao consegui sair do buraco;
3 - como eu vou sair do buraco agora
atencao --> confianca --> dinheiro
papo de cola em leitura em Z
a analise combinatoria apos fazer leitura de baixo para cima
o objetivo do papo que cola fazer a cabeça 3x antes de receber a oferta isso faz gerar autoridade proximidade e aumenta a sua carga de auntenticidade.
entender os passos de construição de um bom papo que cola!!!!!
ALÉM DE FAZER BEM FEITO FAÇA DO SEU JEITO!!!!!!!!!!!!')
$ws.Range("F55").AddComment('Alan Jose Nascimento:
Dica 1 :
Trabalhe com 1 informacção por post-it
Dica 2 :
Sempre coloque um verbo (seu cliente "tem")
Dica 3 :
Nao ofenda seu cliente 
Dica 4 :
nao coloque informações tendenciosas
Dica 5:
Nao coloque informacoes que seja importatnte apenas para voce
Dica 6 :
tire tudo que nao fizer sentido
Dica 7:
Evolua sempre
')
$ws.Range("F57").AddComment('Alan Jose Nascimento:
Pratica é o que leva a perfeição
')
$ws.Range("F66").AddComment('Alan Jose Nascimento:
fazer um flow, o cliente precisa navegar 
atração - experimentação - fechamento (base do humanes)
Replay dos encontrols online 2 e 3
webinar tem detalhes no replay nr 3, trabalhar no funil o que ja deu conta de vender na rua
Mentalidade comunicativa.
Criar um quizz, com o nome  Raio X - Financeiro, qualificar a possivel lead de 3 formas (perfil) .
nota baixa <5   mostrar perigos de continuar assim 
nota media entre 5 e 8 mostrar perigos de continuar assim 
nota alta acima de 8 mostrar perigos de continuar assim 
se esta experimentacao der certo vendeu converteu, fazer webinar 
Ferramentas utilizadas para implementação desse funil:
Mapa Mental: MindMeister ok.
Quiz (teste): OutGrow https://outgrow.co/
Robô (Bot): ManyChat https://manychat.com/
Integrador: Zapier https://zapier.com/
')

$ws.Range("A1").Select() | Out-Null

